$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that receive a plain numeric-looking string must be forced to Text
# format first, otherwise Excel auto-converts them to Number on assignment.
$textCells = @("D5","D6","D7","D8","D10","D11","D12","D13","D16","D19","D20","D21","D23","D24","D25","D26","D28","D30","D32","D34","D35","D36","D38","D39","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "57.941.37"
$ws.Range("E2").Value = "  -1.50%  "
$ws.Range("D3").Value = "2.460.85"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "517.80"
$ws.Range("E5").Value = "  -2.69%  "
$ws.Range("D6").Value = "132.72"
$ws.Range("E6").Value = "  -1.97%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "0.555"
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("D9").Value = "2.466.79"
$ws.Range("E9").Value = "  -1.49%  "
$ws.Range("D10").Value = "0.0974"
$ws.Range("E10").Value = "  -3.42%  "
$ws.Range("D11").Value = "0.157"
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").Value = "5.29"
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("D13").Value = "0.336"
$ws.Range("E13").Value = "  -2.98%  "
$ws.Range("D14").Value = "2.895.09"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("D15").Value = "57.847.33"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").Value = "21.91"
$ws.Range("E16").Value = "  -2.82%  "
$ws.Range("E17").Value = "  -2.55%  "
$ws.Range("D18").Value = "2.459.69"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("D19").Value = "10.56"
$ws.Range("E19").Value = "  -3.52%  "
$ws.Range("D20").Value = "318.42"
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("D21").Value = "4.12"
$ws.Range("E21").Value = "  -2.40%  "
$ws.Range("D23").Value = "5.68"
$ws.Range("E23").Value = "  -4.53%  "
$ws.Range("D24").Value = "64.34"
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("D25").Value = "0.406"
$ws.Range("E25").Value = "  -2.86%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -2.60%  "
$ws.Range("D28").Value = "7.29"
$ws.Range("E28").Value = "  -2.15%  "
$ws.Range("D29").Value = "0.0₃0739"
$ws.Range("E29").Value = "  -2.05%  "
$ws.Range("D30").Value = "168.22"
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("E31").Value = "  -3.00%  "
$ws.Range("D32").Value = "6.20"
$ws.Range("E32").Value = "  -3.00%  "
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "17.97"
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("E37").Value = "  -1.67%  "
$ws.Range("D38").Value = "3.93"
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("D39").Value = "36.35"
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("E40").Value = "  -4.33%  "
$ws.Range("D41").Value = "0.786"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "3.41"
$ws.Range("E42").Value = "  -3.77%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "270.59"
$ws.Range("E43").Value = "  -3.46%  "
$ws.Range("D44").Value = "4.94"
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("D45").Value = "0.588"
$ws.Range("E45").Value = "  -2.48%  "
$ws.Range("D46").Value = "123.38"
$ws.Range("E46").Value = "  -4.57%  "
$ws.Range("D47").Value = "0.0906"
$ws.Range("E47").Value = "  -1.55%  "
$ws.Range("D48").Value = "0.0483"
$ws.Range("E48").Value = "  -2.61%  "
$ws.Range("D49").Value = "0.0211"
$ws.Range("E49").Value = "  -2.92%  "
$ws.Range("D50").Value = "16.72"
$ws.Range("E50").Value = "  -2.13%  "
$ws.Range("D51").Value = "1.723.45"
$ws.Range("E51").Value = "  -1.41%  "
